$d = $word.ActiveDocument

# The document ends with a trailing empty bulleted ("ListParagraph", numId=2)
# paragraph. We add a brand-new list item just before it that reads:
#   "Updating forked repo online on github: <link>"
# leaving the pre-existing trailing empty paragraph untouched.

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# Create a new paragraph right before the trailing empty one; it inherits
# the same paragraph formatting (ListParagraph style / numPr numId=2).
$lastPara.Range.InsertParagraphBefore()

# Re-fetch the freshly inserted (still empty) paragraph - it is now the
# second to last paragraph in the document.
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)

$linkUrl = "https://www.youtube.com/watch?v=TsUIRtT80QU"

# Fill it in as plain text first (label + raw URL) - this cleanly occupies
# the paragraph's placeholder run with no left-over empty runs.
$newPara.Range.Text = "Updating forked repo online on github: " + $linkUrl

# Re-fetch the paragraph range and locate the URL text we just inserted so
# we can convert exactly that span into a hyperlink run, in place.
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$urlRange = $newPara.Range.Duplicate
$found = $urlRange.Find.Execute($linkUrl, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $d.Hyperlinks.Add($urlRange, $linkUrl, [Type]::Missing, [Type]::Missing, $linkUrl) | Out-Null
}
